# Update column F ("人气"/popularity-type numeric) values on the
# "展览" (sheet 1) and "全部类型" (sheet 4) worksheets to match the
# regenerated data output (commit "Update gh-pages to output generated
# at 456a3b4"). Sheets "演出" and "本地生活" are unaffected.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st sheet) - cell -> new value
$updatesSheet1 = @{
    "F4"  = 3488
    "F5"  = 3488
    "F6"  = 246
    "F7"  = 5017
    "F8"  = 506
    "F9"  = 334
    "F11" = 669
    "F13" = 70
    "F14" = 27
    "F15" = 685
    "F16" = 302
    "F22" = 4858
    "F23" = 44
    "F24" = 40
    "F26" = 5974
    "F29" = 3211
    "F30" = 315
    "F31" = 695
    "F32" = 4438
    "F34" = 111
    "F35" = 138
    "F36" = 961
    "F40" = 844
    "F41" = 937
}

# Sheet "全部类型" (4th sheet) - same logical rows, offset by 4 - cell -> new value
$updatesSheet4 = @{
    "F8"  = 3488
    "F9"  = 3488
    "F10" = 246
    "F11" = 5017
    "F12" = 506
    "F13" = 334
    "F15" = 669
    "F16" = 70
    "F17" = 27
    "F18" = 685
    "F19" = 302
    "F26" = 4858
    "F27" = 44
    "F28" = 40
    "F30" = 5974
    "F33" = 3211
    "F34" = 315
    "F35" = 695
    "F36" = 4438
    "F39" = 111
    "F40" = 138
    "F41" = 961
    "F45" = 844
    "F46" = 937
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($cell in $updatesSheet1.Keys) {
    $ws1.Range($cell).Value = $updatesSheet1[$cell]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($cell in $updatesSheet4.Keys) {
    $ws4.Range($cell).Value = $updatesSheet4[$cell]
}
